$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(105, 1).Value = 2016
$ws.Cells.Item(105, 2).Value = 774999
$ws.Cells.Item(105, 3).Value = 1
$ws.Cells.Item(105, 4).Value = "[1596]"
$ws.Cells.Item(105, 5).Value = 1596
$ws.Cells.Item(105, 6).Value = 2003
$ws.Cells.Item(105, 7).Value = 0.7968

$ws.Cells.Item(106, 1).Value = 2020
$ws.Cells.Item(106, 2).Value = 349999
$ws.Cells.Item(106, 3).Value = 1
$ws.Cells.Item(106, 4).Value = "[1644]"
$ws.Cells.Item(106, 5).Value = 1644
$ws.Cells.Item(106, 6).Value = 2044
$ws.Cells.Item(106, 7).Value = 0.8043

$ws.Cells.Item(107, 1).Value = 2024
$ws.Cells.Item(107, 2).Value = 1074999
$ws.Cells.Item(107, 3).Value = 1
$ws.Cells.Item(107, 4).Value = "[1581]"
$ws.Cells.Item(107, 5).Value = 1581
$ws.Cells.Item(107, 6).Value = 1989
$ws.Cells.Item(107, 7).Value = 0.7949000000000001

$ws.Cells.Item(108, 1).Value = 2028
$ws.Cells.Item(108, 2).Value = 274999
$ws.Cells.Item(108, 3).Value = 1
$ws.Cells.Item(108, 4).Value = "[1583]"
$ws.Cells.Item(108, 5).Value = 1583
$ws.Cells.Item(108, 6).Value = 1993
$ws.Cells.Item(108, 7).Value = 0.7943

$ws.Cells.Item(109, 1).Value = 2032
$ws.Cells.Item(109, 2).Value = 574999
$ws.Cells.Item(109, 3).Value = 1
$ws.Cells.Item(109, 4).Value = "[1586]"
$ws.Cells.Item(109, 5).Value = 1586
$ws.Cells.Item(109, 6).Value = 2000
$ws.Cells.Item(109, 7).Value = 0.793

$ws.Cells.Item(110, 1).Value = 2036
$ws.Cells.Item(110, 2).Value = 374999
$ws.Cells.Item(110, 3).Value = 1
$ws.Cells.Item(110, 4).Value = "[1701]"
$ws.Cells.Item(110, 5).Value = 1701
$ws.Cells.Item(110, 6).Value = 2083
$ws.Cells.Item(110, 7).Value = 0.8166

$ws.Cells.Item(111, 1).Value = 2041
$ws.Cells.Item(111, 2).Value = 599999
$ws.Cells.Item(111, 3).Value = 1
$ws.Cells.Item(111, 4).Value = "[1605]"
$ws.Cells.Item(111, 5).Value = 1605
$ws.Cells.Item(111, 6).Value = 1971
$ws.Cells.Item(111, 7).Value = 0.8143

$ws.Cells.Item(112, 1).Value = 2045
$ws.Cells.Item(112, 2).Value = 449999
$ws.Cells.Item(112, 3).Value = 1
$ws.Cells.Item(112, 4).Value = "[1638]"
$ws.Cells.Item(112, 5).Value = 1638
$ws.Cells.Item(112, 6).Value = 2028
$ws.Cells.Item(112, 7).Value = 0.8077

$ws.Cells.Item(113, 1).Value = 2048
$ws.Cells.Item(113, 2).Value = 499999
$ws.Cells.Item(113, 3).Value = 1
$ws.Cells.Item(113, 4).Value = "[1545]"
$ws.Cells.Item(113, 5).Value = 1545
$ws.Cells.Item(113, 6).Value = 1939
$ws.Cells.Item(113, 7).Value = 0.7968

$ws.Cells.Item(114, 1).Value = 2053
$ws.Cells.Item(114, 2).Value = 224999
$ws.Cells.Item(114, 3).Value = 1
$ws.Cells.Item(114, 4).Value = "[1587]"
$ws.Cells.Item(114, 5).Value = 1587
$ws.Cells.Item(114, 6).Value = 1962
$ws.Cells.Item(114, 7).Value = 0.8089

$ws.Cells.Item(115, 1).Value = 2057
$ws.Cells.Item(115, 2).Value = 349999
$ws.Cells.Item(115, 3).Value = 1
$ws.Cells.Item(115, 4).Value = "[1612]"
$ws.Cells.Item(115, 5).Value = 1612
$ws.Cells.Item(115, 6).Value = 2013
$ws.Cells.Item(115, 7).Value = 0.8008

$ws.Cells.Item(116, 1).Value = 2061
$ws.Cells.Item(116, 2).Value = 824999
$ws.Cells.Item(116, 3).Value = 1
$ws.Cells.Item(116, 4).Value = "[1637]"
$ws.Cells.Item(116, 5).Value = 1637
$ws.Cells.Item(116, 6).Value = 2002
$ws.Cells.Item(116, 7).Value = 0.8177

$ws.Cells.Item(117, 1).Value = 2064
$ws.Cells.Item(117, 2).Value = 249999
$ws.Cells.Item(117, 3).Value = 1
$ws.Cells.Item(117, 4).Value = "[1562]"
$ws.Cells.Item(117, 5).Value = 1562
$ws.Cells.Item(117, 6).Value = 1956
$ws.Cells.Item(117, 7).Value = 0.7986

$ws.Cells.Item(118, 1).Value = 2068
$ws.Cells.Item(118, 2).Value = 424999
$ws.Cells.Item(118, 3).Value = 1
$ws.Cells.Item(118, 4).Value = "[1532]"
$ws.Cells.Item(118, 5).Value = 1532
$ws.Cells.Item(118, 6).Value = 1941
$ws.Cells.Item(118, 7).Value = 0.7893

$ws.Cells.Item(119, 1).Value = 2073
$ws.Cells.Item(119, 2).Value = 249999
$ws.Cells.Item(119, 3).Value = 1
$ws.Cells.Item(119, 4).Value = "[1552]"
$ws.Cells.Item(119, 5).Value = 1552
$ws.Cells.Item(119, 6).Value = 2004
$ws.Cells.Item(119, 7).Value = 0.7745

$ws.Cells.Item(120, 1).Value = 2077
$ws.Cells.Item(120, 2).Value = 524999
$ws.Cells.Item(120, 3).Value = 1
$ws.Cells.Item(120, 4).Value = "[1608]"
$ws.Cells.Item(120, 5).Value = 1608
$ws.Cells.Item(120, 6).Value = 2018
$ws.Cells.Item(120, 7).Value = 0.7968

$ws.Cells.Item(121, 1).Value = 2081
$ws.Cells.Item(121, 2).Value = 374999
$ws.Cells.Item(121, 3).Value = 1
$ws.Cells.Item(121, 4).Value = "[1643]"
$ws.Cells.Item(121, 5).Value = 1643
$ws.Cells.Item(121, 6).Value = 2037
$ws.Cells.Item(121, 7).Value = 0.8066

$ws.Cells.Item(122, 1).Value = 2085
$ws.Cells.Item(122, 2).Value = 249999
$ws.Cells.Item(122, 3).Value = 1
$ws.Cells.Item(122, 4).Value = "[1629]"
$ws.Cells.Item(122, 5).Value = 1629
$ws.Cells.Item(122, 6).Value = 2001
$ws.Cells.Item(122, 7).Value = 0.8141

$ws.Cells.Item(123, 1).Value = 2089
$ws.Cells.Item(123, 2).Value = 799999
$ws.Cells.Item(123, 3).Value = 1
$ws.Cells.Item(123, 4).Value = "[1595]"
$ws.Cells.Item(123, 5).Value = 1595
$ws.Cells.Item(123, 6).Value = 1985
$ws.Cells.Item(123, 7).Value = 0.8035

$ws.Cells.Item(124, 1).Value = 2093
$ws.Cells.Item(124, 2).Value = 399999
$ws.Cells.Item(124, 3).Value = 1
$ws.Cells.Item(124, 4).Value = "[1572]"
$ws.Cells.Item(124, 5).Value = 1572
$ws.Cells.Item(124, 6).Value = 1959
$ws.Cells.Item(124, 7).Value = 0.8025

$ws.Cells.Item(125, 1).Value = 2097
$ws.Cells.Item(125, 2).Value = 199999
$ws.Cells.Item(125, 3).Value = 1
$ws.Cells.Item(125, 4).Value = "[1561]"
$ws.Cells.Item(125, 5).Value = 1561
$ws.Cells.Item(125, 6).Value = 1966
$ws.Cells.Item(125, 7).Value = 0.794

$ws.Cells.Item(126, 1).Value = 2101
$ws.Cells.Item(126, 2).Value = 549999
$ws.Cells.Item(126, 3).Value = 1
$ws.Cells.Item(126, 4).Value = "[1625]"
$ws.Cells.Item(126, 5).Value = 1625
$ws.Cells.Item(126, 6).Value = 2024
$ws.Cells.Item(126, 7).Value = 0.8028999999999999

$ws.Cells.Item(127, 1).Value = 2105
$ws.Cells.Item(127, 2).Value = 524999
$ws.Cells.Item(127, 3).Value = 1
$ws.Cells.Item(127, 4).Value = "[1621]"
$ws.Cells.Item(127, 5).Value = 1621
$ws.Cells.Item(127, 6).Value = 2019
$ws.Cells.Item(127, 7).Value = 0.8028999999999999

$ws.Cells.Item(128, 1).Value = 2109
$ws.Cells.Item(128, 2).Value = 674999
$ws.Cells.Item(128, 3).Value = 1
$ws.Cells.Item(128, 4).Value = "[1667]"
$ws.Cells.Item(128, 5).Value = 1667
$ws.Cells.Item(128, 6).Value = 2007
$ws.Cells.Item(128, 7).Value = 0.8306

$ws.Cells.Item(129, 1).Value = 2112
$ws.Cells.Item(129, 2).Value = 49999
$ws.Cells.Item(129, 3).Value = 1
$ws.Cells.Item(129, 4).Value = "[1621]"
$ws.Cells.Item(129, 5).Value = 1621
$ws.Cells.Item(129, 6).Value = 1994
$ws.Cells.Item(129, 7).Value = 0.8129

$ws.Cells.Item(130, 1).Value = 2116
$ws.Cells.Item(130, 2).Value = 174999
$ws.Cells.Item(130, 3).Value = 1
$ws.Cells.Item(130, 4).Value = "[1626]"
$ws.Cells.Item(130, 5).Value = 1626
$ws.Cells.Item(130, 6).Value = 2006
$ws.Cells.Item(130, 7).Value = 0.8106

$ws.Cells.Item(131, 1).Value = 2120
$ws.Cells.Item(131, 2).Value = 1149999
$ws.Cells.Item(131, 3).Value = 1
$ws.Cells.Item(131, 4).Value = "[1614]"
$ws.Cells.Item(131, 5).Value = 1614
$ws.Cells.Item(131, 6).Value = 1988
$ws.Cells.Item(131, 7).Value = 0.8119

$ws.Cells.Item(132, 1).Value = 2124
$ws.Cells.Item(132, 2).Value = 524999
$ws.Cells.Item(132, 3).Value = 1
$ws.Cells.Item(132, 4).Value = "[1618]"
$ws.Cells.Item(132, 5).Value = 1618
$ws.Cells.Item(132, 6).Value = 1998
$ws.Cells.Item(132, 7).Value = 0.8098

$ws.Cells.Item(133, 1).Value = 2128
$ws.Cells.Item(133, 2).Value = 724999
$ws.Cells.Item(133, 3).Value = 1
$ws.Cells.Item(133, 4).Value = "[1613]"
$ws.Cells.Item(133, 5).Value = 1613
$ws.Cells.Item(133, 6).Value = 1999
$ws.Cells.Item(133, 7).Value = 0.8069

$ws.Cells.Item(134, 1).Value = 2132
$ws.Cells.Item(134, 2).Value = 1024999
$ws.Cells.Item(134, 3).Value = 1
$ws.Cells.Item(134, 4).Value = "[1682]"
$ws.Cells.Item(134, 5).Value = 1682
$ws.Cells.Item(134, 6).Value = 2070
$ws.Cells.Item(134, 7).Value = 0.8126

$ws.Cells.Item(135, 1).Value = 2136
$ws.Cells.Item(135, 2).Value = 249999
$ws.Cells.Item(135, 3).Value = 1
$ws.Cells.Item(135, 4).Value = "[1602]"
$ws.Cells.Item(135, 5).Value = 1602
$ws.Cells.Item(135, 6).Value = 1990
$ws.Cells.Item(135, 7).Value = 0.805

$ws.Cells.Item(136, 1).Value = 2140
$ws.Cells.Item(136, 2).Value = 474999
$ws.Cells.Item(136, 3).Value = 1
$ws.Cells.Item(136, 4).Value = "[1630]"
$ws.Cells.Item(136, 5).Value = 1630
$ws.Cells.Item(136, 6).Value = 1997
$ws.Cells.Item(136, 7).Value = 0.8162
